$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the candidate database (A1:K1)
$ws.Range("A1").Value = "Título da vaga desejada"
$ws.Range("B1").Value = "Tipo da vaga desejada"
$ws.Range("C1").Value = "Área de interesse"
$ws.Range("D1").Value = "Nível de inglês"
$ws.Range("E1").Value = "Nível de espanhol"
$ws.Range("F1").Value = "Outros idiomas"
$ws.Range("G1").Value = "Competências técnicas"
$ws.Range("H1").Value = "Competências comportamentais"
$ws.Range("I1").Value = "Disponível para viagens? (Sim/Não)"
$ws.Range("J1").Value = "Possui equipamento próprio? (Sim/Não)"
$ws.Range("K1").Value = "Expectativa salarial"

# Resize every column to (best) fit the new, longer header text
$ws.Columns.Item(1).ColumnWidth = 21.09
$ws.Columns.Item(2).ColumnWidth = 19.59
$ws.Columns.Item(3).ColumnWidth = 15.75
$ws.Columns.Item(4).ColumnWidth = 13.42
$ws.Columns.Item(5).ColumnWidth = 16.42
$ws.Columns.Item(6).ColumnWidth = 13.59
$ws.Columns.Item(7).ColumnWidth = 20.59
$ws.Columns.Item(8).ColumnWidth = 29.42
$ws.Columns.Item(9).ColumnWidth = 32.25
$ws.Columns.Item(10).ColumnWidth = 36.59
$ws.Columns.Item(11).ColumnWidth = 17.25

# Move the active selection to L11 (matches the final saved cursor position)
$ws.Range("L11").Select() | Out-Null

$wb.Save()
